$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Analise de Eventos")

# --- Row 17: new scenario (event 15) ---
$ws.Range("B17").Value = "FA"
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = "Solicita cancelamento à fábrica"
$ws.Range("F17").Value = "x"

# --- Row 18: new scenario (event 16) ---
$ws.Range("C18").Value = 16
$ws.Range("D18").Value = "Cancela a Nota Fiscal do pedido"

# --- Row 19: new scenario (event 17) ---
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = "Cancelar Pedido"

# --- Remaining shared-string cells (order matches original authoring) ---
$ws.Range("H18").Value = "x(15)"
$ws.Range("G19").Value = "x(16)"

# --- Remove trailing hidden empty rows 21-23 ---
$ws.Rows("21:23").Delete()

# --- Update view / selection state ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("G19").Select()
